# Auto - Update data with bot!
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9: title + link
$ws.Range("D9").Value = "공대가 좋아하는 DS자소서 vs. 공대를 싫어하게 된 DS가 보는 자소서"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/engineers-loved-ds-hated-ds/#utm_source=rss&utm_medium=rss&utm_campaign=engineers-loved-ds-hated-ds"

# Row 28: title + link
$ws.Range("D28").Value = "[GNN 2] Graph LSTM"
$ws.Range("E28").Value = "https://ropiens.tistory.com/195"

# Row 37: title only
$ws.Range("D37").Value = "dsba_seminar"

# Row 52: title + link
$ws.Range("D52").Value = "효율적인 R 프로그래밍"
$ws.Range("E52").Value = "http://ds.sumeun.org/?p=2614"
